$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Affiliate Disclosure: Required in introduction*") {
        $p.Range.Delete()
        break
    }
}
